$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 25000
$ws.Range("I19").Value = 25000
$ws.Range("K19").Value = 25000
$ws.Range("M19").Value = -24825
$ws.Range("H46").Value = 800
$ws.Range("J46").Value = 800
$ws.Range("L46").Value = 2400
$ws.Range("N46").Value = -2638
$ws.Range("H60").Value = 800
$ws.Range("J60").Value = 800
$ws.Range("L60").Value = 2400
$ws.Range("N60").Value = -3368
$ws.Range("H74").Value = 10487.5
$ws.Range("J74").Value = 10487.5
$ws.Range("L74").Value = 10487.5
$ws.Range("N74").Value = -12359.5
$ws.Range("H77").Value = 10487.5
$ws.Range("J77").Value = 10487.5
$ws.Range("L77").Value = 52437.5
$ws.Range("N77").Value = -61797.5
$ws.Range("H137").Value = 2666.6667
$ws.Range("I137").Value = 1750
$ws.Range("J137").Value = 4500
$ws.Range("K137").Value = 5250
$ws.Range("L137").Value = 13500
$ws.Range("M137").Value = -2700
$ws.Range("N137").Value = -18600
$ws.Range("H138").Value = 3361.4375
$ws.Range("J138").Value = 3499
$ws.Range("L138").Value = 10497
$ws.Range("N138").Value = -20777

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3061.8
$ws.Range("I32").Value = 3061.8
$ws.Range("K32").Value = 3061.8
$ws.Range("M32").Value = -2774.8
$ws.Range("H102").Value = 1352.5
$ws.Range("I102").Value = 1499.6666
$ws.Range("J102").Value = 911
$ws.Range("K102").Value = 1499.6666
$ws.Range("L102").Value = 911
$ws.Range("M102").Value = 122.3334
$ws.Range("N102").Value = -4155

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H35").Value = 49920
$ws.Range("J35").Value = 49920
$ws.Range("L35").Value = 49920
$ws.Range("N35").Value = -50540
$ws.Range("H129").Value = 58931
$ws.Range("I129").Value = 58931
$ws.Range("J129").Value = 58931
$ws.Range("K129").Value = 58931
$ws.Range("L129").Value = 58931
$ws.Range("M129").Value = -53931
$ws.Range("N129").Value = -68931
$ws.Range("H134").Value = 2781.2144
$ws.Range("I134").Value = 1419.7142
$ws.Range("J134").Value = 4142.7144
$ws.Range("K134").Value = 4259.142599999999
$ws.Range("L134").Value = 12428.1432
$ws.Range("M134").Value = -1724.142599999999
$ws.Range("N134").Value = -17498.1432

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 47.88889
$ws.Range("I7").Value = 49.57143
$ws.Range("K7").Value = 49.57143
$ws.Range("M7").Value = 63.42857
$ws.Range("H31").Value = 2429.1667
$ws.Range("I31").Value = 1143.75
$ws.Range("K31").Value = 1143.75
$ws.Range("M31").Value = -848.75
$ws.Range("H34").Value = 2429.1667
$ws.Range("I34").Value = 1143.75
$ws.Range("K34").Value = 1143.75
$ws.Range("M34").Value = -941.75
$ws.Range("H86").Value = 5468.1577
$ws.Range("I86").Value = 2610.9
$ws.Range("J86").Value = 8642.888999999999
$ws.Range("K86").Value = 2610.9
$ws.Range("L86").Value = 8642.888999999999
$ws.Range("M86").Value = -1487.9
$ws.Range("N86").Value = -10888.889
$ws.Range("H89").Value = 5468.1577
$ws.Range("I89").Value = 2610.9
$ws.Range("J89").Value = 8642.888999999999
$ws.Range("K89").Value = 13054.5
$ws.Range("L89").Value = 43214.44499999999
$ws.Range("M89").Value = -7438.5
$ws.Range("N89").Value = -54446.44499999999
$ws.Range("H99").Value = 2533.3333
$ws.Range("I99").Value = 2400
$ws.Range("K99").Value = 2400
$ws.Range("M99").Value = -902
$ws.Range("H126").Value = 2533.3333
$ws.Range("I126").Value = 2400
$ws.Range("K126").Value = 7200
$ws.Range("M126").Value = -4730
$ws.Range("H132").Value = 3560.8
$ws.Range("I132").Value = 2944.5715
$ws.Range("K132").Value = 8833.7145
$ws.Range("M132").Value = -6303.7145
$ws.Range("H134").Value = 800
$ws.Range("I134").Value = 800
$ws.Range("K134").Value = 2400
$ws.Range("M134").Value = 135
$ws.Range("H141").Value = 49980.5
$ws.Range("J141").Value = 49980.5
$ws.Range("L141").Value = 49980.5
$ws.Range("N141").Value = -60340.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H6").Value = 30.2
$ws.Range("I6").Value = 37.25
$ws.Range("J6").Value = 2
$ws.Range("K6").Value = 111.75
$ws.Range("L6").Value = 6
$ws.Range("M6").Value = 1.25
$ws.Range("N6").Value = -232
$ws.Range("H12").Value = 57
$ws.Range("I12").Value = 54
$ws.Range("K12").Value = 162
$ws.Range("M12").Value = 11
$ws.Range("H80").Value = 1950.25
$ws.Range("I80").Value = 1966
$ws.Range("K80").Value = 5898
$ws.Range("M80").Value = -4962
$ws.Range("H83").Value = 1950.25
$ws.Range("I83").Value = 1966
$ws.Range("K83").Value = 17694
$ws.Range("M83").Value = -13014
$ws.Range("H128").Value = 129962.664
$ws.Range("I128").Value = 129962.664
$ws.Range("K128").Value = 389887.992
$ws.Range("M128").Value = -384907.992
$ws.Range("H131").Value = 2787.0645
$ws.Range("I131").Value = 1470
$ws.Range("J131").Value = 2928.1785
$ws.Range("K131").Value = 4410
$ws.Range("L131").Value = 8784.5355
$ws.Range("M131").Value = 630
$ws.Range("N131").Value = -18864.5355

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 3338.182
$ws.Range("I102").Value = 2774.5715
$ws.Range("K102").Value = 2774.5715
$ws.Range("M102").Value = -1152.5715

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3999
$ws.Range("I7").Value = 3999
$ws.Range("K7").Value = 3999
$ws.Range("M7").Value = -3887
$ws.Range("H68").Value = 2933.6667
$ws.Range("I68").Value = 1899
$ws.Range("K68").Value = 1899
$ws.Range("M68").Value = -1150
$ws.Range("H71").Value = 2933.6667
$ws.Range("I71").Value = 1899
$ws.Range("K71").Value = 9495
$ws.Range("M71").Value = -5751
$ws.Range("H82").Value = 23789.3
$ws.Range("I82").Value = 19148.166
$ws.Range("J82").Value = 30751
$ws.Range("K82").Value = 19148.166
$ws.Range("L82").Value = 30751
$ws.Range("M82").Value = -18787.166
$ws.Range("N82").Value = -31473
$ws.Range("H85").Value = 23789.3
$ws.Range("I85").Value = 19148.166
$ws.Range("J85").Value = 30751
$ws.Range("K85").Value = 19148.166
$ws.Range("L85").Value = 30751
$ws.Range("M85").Value = -17900.166
$ws.Range("N85").Value = -33247
$ws.Range("H93").Value = 2390
$ws.Range("I93").Value = 2679.6
$ws.Range("J93").Value = 1666
$ws.Range("K93").Value = 2679.6
$ws.Range("L93").Value = 1666
$ws.Range("M93").Value = -1431.6
$ws.Range("N93").Value = -4162
$ws.Range("H126").Value = 3999
$ws.Range("I126").Value = 3999
$ws.Range("K126").Value = 11997
$ws.Range("M126").Value = -9527
$ws.Range("H136").Value = 5500
$ws.Range("I136").Value = 5500
$ws.Range("K136").Value = 16500
$ws.Range("M136").Value = -13950

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 583.7143
$ws.Range("I122").Value = 501.75
$ws.Range("K122").Value = 1505.25
$ws.Range("M122").Value = 944.75
$ws.Range("H126").Value = 1222
$ws.Range("I126").Value = 867.2222
$ws.Range("K126").Value = 2601.6666
$ws.Range("M126").Value = -131.6666
